$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.164.84"
$ws.Range("E2").Value = "  +3.95%  "
$ws.Range("D3").Value = "1.603.50"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.93"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  +3.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0615"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.03"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +4.93%  "
$ws.Range("D12").Value = "1.825.55"
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").Value = "1.598.51"
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "26.147.83"
$ws.Range("E16").Value = "  +3.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.53"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "205.68"
$ws.Range("E20").Value = "  +11.89%  "
$ws.Range("E21").Value = "  +4.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.00"
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  +9.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.88"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.116.28"
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("E37").Value = "  +10.25%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.782"
$ws.Range("E40").Value = "  +4.11%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "1.738.14"
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.96"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  +6.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.40"
$ws.Range("E47").Value = "  +2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0504"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.18"
$ws.Range("E51").Value = "  +0.87%  "
